$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert two new price rows (newest entries) just above the
# existing "1a (guarda)/2a (guarda)" pair that used to sit at rows 70-71,
# pushing the rest of the table (old rows 70-162) down by two rows.
$ws.Rows("70:71").Insert()

# New row 70: "1a nueva(o)"
$ws.Range("A70").Value = 11
$ws.Range("B70").Value = "Vega Monumental Concepción"
$ws.Range("C70").Value = "Bíobío"
$ws.Range("D70").Value = 44579
$ws.Range("E70").Value = 8
$ws.Range("F70").Value = 100112045
$ws.Range("G70").Value = "Zapallo"
$ws.Range("H70").Value = "Camote"
$ws.Range("I70").Value = "1a nueva(o)"
$ws.Range("J70").Value = 800
$ws.Range("K70").Value = 500
$ws.Range("L70").Value = 550
$ws.Range("M70").Value = 525
$ws.Range("N70").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O70").Value = "Región de O'Higgins"
$ws.Range("P70").Value = 525
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"

# New row 71: "2a nueva(o)"
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = "Vega Monumental Concepción"
$ws.Range("C71").Value = "Bíobío"
$ws.Range("D71").Value = 44579
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 100112045
$ws.Range("G71").Value = "Zapallo"
$ws.Range("H71").Value = "Camote"
$ws.Range("I71").Value = "2a nueva(o)"
$ws.Range("J71").Value = 400
$ws.Range("K71").Value = 450
$ws.Range("L71").Value = 450
$ws.Range("M71").Value = 450
$ws.Range("N71").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O71").Value = "Región de O'Higgins"
$ws.Range("P71").Value = 450
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = "Hortaliza"
